$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}

Replace-Text "2025-01-08 Wednesday" "2025-01-14 Tuesday"

Replace-Text "802×6=" "867×8="
Replace-Text "103×2=" "222×4="
Replace-Text "330×9=" "364×5="
Replace-Text "679×8=" "598×4="
Replace-Text "566×5=" "850×9="
Replace-Text "736×6=" "651×5="
Replace-Text "622×5=" "603×7="
Replace-Text "560×4=" "203×6="
Replace-Text "730×5=" "623×5="
Replace-Text "716×9=" "336×6="
Replace-Text "140×8=" "977×7="
Replace-Text "827×2=" "593×5="
Replace-Text "641×4=" "166×5="
Replace-Text "463×7=" "936×4="
Replace-Text "981×6=" "768×5="
Replace-Text "900×6=" "869×5="
Replace-Text "407×9=" "118×9="
Replace-Text "955×6=" "539×8="
Replace-Text "581×6=" "988×4="
Replace-Text "359×4=" "781×8="
Replace-Text "824×2=" "279×9="
Replace-Text "931×5=" "766×5="
Replace-Text "298×5=" "250×3="
Replace-Text "214×4=" "610×9="
Replace-Text "892×9=" "501×8="
